$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy A26:B26 (U7 / Power Switch) -> C26:D26
$ws.Range("A26:B26").Copy()
$ws.Range("C26").PasteSpecial(-4104)  # xlPasteAll
$ws.Range("A26").Value = "MIC94093YC6"
$ws.Range("B26").Value = "SC-70-6"
